$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume/change (column E) values
# to reflect the latest scrape, per GitHub Actions update job.
# Values that parse as plain numbers (single decimal point) are prefixed
# with a leading apostrophe so Excel stores them as text, matching the
# original inline-string cell formatting (e.g. "491.81", not 491.81).
$ws.Range('D2').Value = '56.603.67'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '2.488.20'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''491.81'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').Value = '''152.45'
$ws.Range('E6').Value = '  +7.09%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  -0.93%  '
$ws.Range('D9').Value = '2.501.94'
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('E10').Value = '  +3.85%  '
$ws.Range('D11').Value = '''0.0984'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('D12').Value = '''0.334'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').Value = '2.921.83'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').Value = '56.830.31'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').Value = '''21.32'
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').Value = '2.497.67'
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('D19').Value = '''4.55'
$ws.Range('E19').Value = '  +2.70%  '
$ws.Range('D20').Value = '''10.33'
$ws.Range('E20').Value = '  +2.24%  '
$ws.Range('D21').Value = '''320.58'
$ws.Range('E21').Value = '  -1.32%  '
$ws.Range('D22').Value = '''0.998'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '''5.89'
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('D24').Value = '''58.68'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.86%  '
$ws.Range('D27').Value = '''0.162'
$ws.Range('E27').Value = '  -5.15%  '
$ws.Range('D28').Value = '2.602.82'
$ws.Range('E28').Value = '  -1.02%  '
$ws.Range('D29').Value = '''7.59'
$ws.Range('E29').Value = '  +0.71%  '
$ws.Range('D30').Value = '0.0₃0806'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('D32').Value = '''150.75'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').Value = '''18.35'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('E35').Value = '  +0.69%  '
$ws.Range('E36').Value = '  +3.04%  '
$ws.Range('D37').Value = '''3.78'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('D38').Value = '''0.869'
$ws.Range('E38').Value = '  -3.20%  '
$ws.Range('D39').Value = '''1.39'
$ws.Range('E39').Value = '  +3.94%  '
$ws.Range('E40').Value = '  -0.86%  '
$ws.Range('E41').Value = '  +1.95%  '
$ws.Range('D42').Value = '''0.0564'
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('D43').Value = '''0.616'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').Value = '''4.89'
$ws.Range('E45').Value = '  +2.09%  '
$ws.Range('D46').Value = '''268.66'
$ws.Range('E46').Value = '  +3.60%  '
$ws.Range('D47').Value = '''0.0930'
$ws.Range('E47').Value = '  +1.49%  '
$ws.Range('D48').Value = '''0.0229'
$ws.Range('E48').Value = '  +0.84%  '
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').Value = '''17.77'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').Value = '1.888.71'
$ws.Range('E51').Value = '  -6.41%  '